# Updates the Fixed Residue PSSM matrix (B2:K21) with supplemental-figure recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-18.71168377587409,2.448340879499861,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    3 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,2.379152937443538,-18.71168377587409,-18.71168377587409)
    4 = @(-18.71168377587409,2.179982021358845,2.867172129873833,-18.71168377587409,2.410158928815775,-18.71168377587409,1.914239950257582,-18.71168377587409,2.398950545940953,-18.71168377587409)
    5 = @(-18.71168377587409,0.9703043574857398,-18.71168377587409,-18.71168377587409,-18.71168377587409,2.064528348828615,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    6 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    7 = @(2.97923709316749,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    8 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,2.867991704071978,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    9 = @(3.59863527393555,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    10 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,1.483353144772655,-18.71168377587409,-18.71168377587409)
    11 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,2.032786300926811,-18.71168377587409,2.760239955781004,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    12 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    13 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,1.672839671139044,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,2.469488503721572,4.32192490246357)
    14 = @(-18.71168377587409,-18.71168377587409,1.626444320204386,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    15 = @(-18.71168377587409,-18.71168377587409,-0.2484428032783558,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409)
    16 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,2.233140505287011,-18.71168377587409)
    17 = @(-18.71168377587409,0.6414738208519519,-0.06093702419174819,-18.71168377587409,-18.71168377587409,-18.71168377587409,0.4166454147062919,0.8620983769292605,0.9192960584204212,-18.71168377587409)
    18 = @(-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,-18.71168377587409,0.3777072193286422,0.7712206628980539,1.375479905948339,-18.71168377587409)
    19 = @(-18.71168377587409,-18.71168377587409,1.723946812030946,-18.71168377587409,-18.71168377587409,-18.71168377587409,1.887568955831113,2.189428386589786,-18.71168377587409,-18.71168377587409)
    20 = @(-18.71168377587409,1.661093093612661,2.17385844623418,-18.71168377587409,3.876234942140015,-18.71168377587409,2.171040610243875,1.969394175756064,-18.71168377587409,-18.71168377587409)
    21 = @(-18.71168377587409,1.735033054156498,-18.71168377587409,2.43805900287885,-18.71168377587409,3.17662582000004,2.43112194257548,-18.71168377587409,-18.71168377587409,-18.71168377587409)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $rowVals[$i]
    }
}
